# Update StructureDefinition-communication-product workbook:
# - Metadata sheet: URL, Version, Date, Publisher
# - Elements sheet: clear stray Constraint(s) text on Extension row,
#   and update embedded URLs (Fixed Value + Binding Value Set)

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/communication-product"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/communication-product"
$elements.Range("Y7").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/eng-product"

# Column Y (Binding Value Set) grows to fit the longer URL - match the
# autofit width the desktop app would have produced for the new text.
$elements.Columns.Item(25).ColumnWidth = 51.166666666666664
